$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 9-12), mirroring the existing structure (A: date serial, B..M numeric, N: "Named")
$rows = @(
    @(42613.758136574077, -6,  48, 50, 48, 56, 17612, 13412, 748, 105, 110, 7, 9),
    @(42613.88585648148,  -10, 47, 51, 47, 56, 15606, 14787, 823, 110, 119, 7, 9),
    @(42614.884236111109, -20, 49, 49, 49, 88, 13695, 10342, 565, 91,  92,  1, 8),
    @(42615.884722222225, -16, 51, 48, 51, 99, 14357, 12511, 612, 115, 107, 0, 7)
)

$startRow = 9
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    for ($c = 2; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
    $ws.Cells.Item($r, 14).Value = "Named"
}
